$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(306, 1).Value = 304
$ws.Cells.Item(306, 2).Value = 'Aldnoah Zero 2'
$ws.Cells.Item(306, 3).Value = 3
$ws.Cells.Item(306, 4).Value = 1
$ws.Cells.Item(306, 5).Value = 'Winter 2015'

$ws.Cells.Item(307, 1).Value = 305
$ws.Cells.Item(307, 2).Value = 'Absolute Duo'
$ws.Cells.Item(307, 3).Value = 1
$ws.Cells.Item(307, 4).Value = 3
$ws.Cells.Item(307, 5).Value = 'Winter 2015'

$ws.Cells.Item(308, 1).Value = 306
$ws.Cells.Item(308, 2).Value = 'Dog Days" (Season 3)'
$ws.Cells.Item(308, 3).Value = 1
$ws.Cells.Item(308, 4).Value = 3
$ws.Cells.Item(308, 5).Value = 'Winter 2015'

$ws.Cells.Item(309, 1).Value = 307
$ws.Cells.Item(309, 2).Value = 'Durarara!!x2 Shou'
$ws.Cells.Item(309, 3).Value = 1
$ws.Cells.Item(309, 4).Value = 3
$ws.Cells.Item(309, 5).Value = 'Winter 2015'

$ws.Cells.Item(310, 1).Value = 308
$ws.Cells.Item(310, 2).Value = 'Isuca'
$ws.Cells.Item(310, 3).Value = 1
$ws.Cells.Item(310, 4).Value = 3
$ws.Cells.Item(310, 5).Value = 'Winter 2015'

$ws.Cells.Item(311, 1).Value = 309
$ws.Cells.Item(311, 2).Value = 'Juuou Mujin no Fafnir'
$ws.Cells.Item(311, 3).Value = 1
$ws.Cells.Item(311, 4).Value = 3
$ws.Cells.Item(311, 5).Value = 'Winter 2015'

$ws.Cells.Item(312, 1).Value = 310
$ws.Cells.Item(312, 2).Value = 'Kamisama Hajimemashita 2nd'
$ws.Cells.Item(312, 3).Value = 1
$ws.Cells.Item(312, 4).Value = 3
$ws.Cells.Item(312, 5).Value = 'Winter 2015'

$ws.Cells.Item(313, 1).Value = 311
$ws.Cells.Item(313, 2).Value = 'Kantai Collection : Kan Colle'
$ws.Cells.Item(313, 3).Value = 1
$ws.Cells.Item(313, 4).Value = 3
$ws.Cells.Item(313, 5).Value = 'Winter 2015'

$ws.Cells.Item(314, 1).Value = 312
$ws.Cells.Item(314, 2).Value = 'Koufuku Graffiti [SHAFT]'
$ws.Cells.Item(314, 3).Value = 1
$ws.Cells.Item(314, 4).Value = 3
$ws.Cells.Item(314, 5).Value = 'Winter 2015'

$ws.Cells.Item(315, 1).Value = 313
$ws.Cells.Item(315, 2).Value = 'Saenai Heroine no Sodate-kata'
$ws.Cells.Item(315, 3).Value = 1
$ws.Cells.Item(315, 4).Value = 3
$ws.Cells.Item(315, 5).Value = 'Winter 2015'

$ws.Cells.Item(316, 1).Value = 314
$ws.Cells.Item(316, 2).Value = 'THE iDOLM@STER : Cinderella Girls'
$ws.Cells.Item(316, 3).Value = 1
$ws.Cells.Item(316, 4).Value = 3
$ws.Cells.Item(316, 5).Value = 'Winter 2015'

$ws.Cells.Item(317, 1).Value = 315
$ws.Cells.Item(317, 2).Value = 'Yuri Kuma Arashi'
$ws.Cells.Item(317, 3).Value = 1
$ws.Cells.Item(317, 4).Value = 3
$ws.Cells.Item(317, 5).Value = 'Winter 2015'

$ws.Cells.Item(318, 1).Value = 316
$ws.Cells.Item(318, 2).Value = 'Hori-san to Miyamura-kun: Shingakki'
$ws.Cells.Item(318, 3).Value = 1
$ws.Cells.Item(318, 4).Value = 2
$ws.Cells.Item(318, 5).Value = 'Winter 2015'

$ws.Cells.Item(319, 1).Value = 317
$ws.Cells.Item(319, 2).Value = 'Shimai Maou no Keiyakusha'
$ws.Cells.Item(319, 3).Value = 1
$ws.Cells.Item(319, 4).Value = 2
$ws.Cells.Item(319, 5).Value = 'Winter 2015'

$ws.Cells.Item(320, 1).Value = 318
$ws.Cells.Item(320, 2).Value = 'Ansatsu Kyoushitsu'
$ws.Cells.Item(320, 3).Value = 1
$ws.Cells.Item(320, 4).Value = 1
$ws.Cells.Item(320, 5).Value = 'Winter 2015'

$ws.Cells.Item(321, 1).Value = 319
$ws.Cells.Item(321, 2).Value = 'Junketsu no Maria'
$ws.Cells.Item(321, 3).Value = 3
$ws.Cells.Item(321, 4).Value = 1
$ws.Cells.Item(321, 5).Value = 'Winter 2015'

$ws.Cells.Item(322, 1).Value = 320
$ws.Cells.Item(322, 2).Value = 'Rolling Girls'
$ws.Cells.Item(322, 3).Value = 1
$ws.Cells.Item(322, 4).Value = 1
$ws.Cells.Item(322, 5).Value = 'Winter 2015'

$ws.Cells.Item(323, 1).Value = 321
$ws.Cells.Item(323, 2).Value = 'Seiken Tsukai no World Break'
$ws.Cells.Item(323, 3).Value = 1
$ws.Cells.Item(323, 4).Value = 1
$ws.Cells.Item(323, 5).Value = 'Winter 2015'

$ws.Cells.Item(324, 1).Value = 322
$ws.Cells.Item(324, 2).Value = 'Tantei Opera Milky Holmes TD'
$ws.Cells.Item(324, 3).Value = 1
$ws.Cells.Item(324, 4).Value = 1
$ws.Cells.Item(324, 5).Value = 'Winter 2015'

$ws.Cells.Item(325, 1).Value = 323
$ws.Cells.Item(325, 2).Value = 'Tokyo Ghoul VA'
$ws.Cells.Item(325, 3).Value = 2
$ws.Cells.Item(325, 4).Value = 1
$ws.Cells.Item(325, 5).Value = 'Winter 2015'

$ws.Cells.Item(326, 1).Value = 324
$ws.Cells.Item(326, 2).Value = 'Nisekoi 2'
$ws.Cells.Item(326, 3).Value = 3
$ws.Cells.Item(326, 4).Value = 3
$ws.Cells.Item(326, 5).Value = 'Spring 2015'

$ws.Cells.Item(327, 1).Value = 325
$ws.Cells.Item(327, 2).Value = 'Yamadakun to Nananin no Majo'
$ws.Cells.Item(327, 3).Value = 3
$ws.Cells.Item(327, 4).Value = 3
$ws.Cells.Item(327, 5).Value = 'Spring 2015'

$ws.Cells.Item(328, 1).Value = 326
$ws.Cells.Item(328, 2).Value = 'Shokugeki no Souma'
$ws.Cells.Item(328, 3).Value = 3
$ws.Cells.Item(328, 4).Value = 2
$ws.Cells.Item(328, 5).Value = 'Spring 2015'

$ws.Cells.Item(329, 1).Value = 327
$ws.Cells.Item(329, 2).Value = 'Dungeon ni Deai o Motomeru no wa Machigatteiru Darou ka?'
$ws.Cells.Item(329, 3).Value = 3
$ws.Cells.Item(329, 4).Value = 2
$ws.Cells.Item(329, 5).Value = 'Spring 2015'

$ws.Cells.Item(330, 1).Value = 328
$ws.Cells.Item(330, 2).Value = 'Denpa Kyoushi'
$ws.Cells.Item(330, 3).Value = 1
$ws.Cells.Item(330, 4).Value = 3
$ws.Cells.Item(330, 5).Value = 'Spring 2015'

$ws.Cells.Item(331, 1).Value = 329
$ws.Cells.Item(331, 2).Value = 'Etotama'
$ws.Cells.Item(331, 3).Value = 1
$ws.Cells.Item(331, 4).Value = 3
$ws.Cells.Item(331, 5).Value = 'Spring 2015'

$ws.Cells.Item(332, 1).Value = 330
$ws.Cells.Item(332, 2).Value = 'Grisaia no Meikyuu'
$ws.Cells.Item(332, 3).Value = 1
$ws.Cells.Item(332, 4).Value = 3
$ws.Cells.Item(332, 5).Value = 'Spring 2015'

$ws.Cells.Item(333, 1).Value = 331
$ws.Cells.Item(333, 2).Value = 'Hello!! Kiniro Mosaic'
$ws.Cells.Item(333, 3).Value = 1
$ws.Cells.Item(333, 4).Value = 3
$ws.Cells.Item(333, 5).Value = 'Spring 2015'

$ws.Cells.Item(334, 1).Value = 332
$ws.Cells.Item(334, 2).Value = 'Hibike! Euphonium'
$ws.Cells.Item(334, 3).Value = 3
$ws.Cells.Item(334, 4).Value = 3
$ws.Cells.Item(334, 5).Value = 'Spring 2015'

$ws.Cells.Item(335, 1).Value = 333
$ws.Cells.Item(335, 2).Value = 'High School DxD BorN'
$ws.Cells.Item(335, 3).Value = 3
$ws.Cells.Item(335, 4).Value = 3
$ws.Cells.Item(335, 5).Value = 'Spring 2015'

$ws.Cells.Item(336, 1).Value = 334
$ws.Cells.Item(336, 2).Value = 'Kyoukai no Rinne'
$ws.Cells.Item(336, 3).Value = 1
$ws.Cells.Item(336, 4).Value = 3
$ws.Cells.Item(336, 5).Value = 'Spring 2015'

$ws.Cells.Item(337, 1).Value = 335
$ws.Cells.Item(337, 2).Value = 'Plastic Memories'
$ws.Cells.Item(337, 3).Value = 3
$ws.Cells.Item(337, 4).Value = 3
$ws.Cells.Item(337, 5).Value = 'Spring 2015'

$ws.Cells.Item(338, 1).Value = 336
$ws.Cells.Item(338, 2).Value = 'Re-Kan!'
$ws.Cells.Item(338, 3).Value = 1
$ws.Cells.Item(338, 4).Value = 3
$ws.Cells.Item(338, 5).Value = 'Spring 2015'

$ws.Cells.Item(339, 1).Value = 337
$ws.Cells.Item(339, 2).Value = 'Urawa no Usagi-chan'
$ws.Cells.Item(339, 3).Value = 1
$ws.Cells.Item(339, 4).Value = 3
$ws.Cells.Item(339, 5).Value = 'Spring 2015'

$ws.Cells.Item(340, 1).Value = 338
$ws.Cells.Item(340, 2).Value = 'Yahari Ore no Seishun Love Comedy wa Machigatteiru. Zoku'
$ws.Cells.Item(340, 3).Value = 1
$ws.Cells.Item(340, 4).Value = 3
$ws.Cells.Item(340, 5).Value = 'Spring 2015'

$ws.Cells.Item(341, 1).Value = 339
$ws.Cells.Item(341, 2).Value = 'Gintama (2015)'
$ws.Cells.Item(341, 3).Value = 1
$ws.Cells.Item(341, 4).Value = 2
$ws.Cells.Item(341, 5).Value = 'Spring 2015'

$ws.Cells.Item(342, 1).Value = 340
$ws.Cells.Item(342, 2).Value = 'Mahou Shoujo Lyrical Nanoha ViVid'
$ws.Cells.Item(342, 3).Value = 1
$ws.Cells.Item(342, 4).Value = 2
$ws.Cells.Item(342, 5).Value = 'Spring 2015'

$ws.Cells.Item(343, 1).Value = 341
$ws.Cells.Item(343, 2).Value = 'Ore Monogatari!!'
$ws.Cells.Item(343, 3).Value = 1
$ws.Cells.Item(343, 4).Value = 2
$ws.Cells.Item(343, 5).Value = 'Spring 2015'

$ws.Cells.Item(344, 1).Value = 342
$ws.Cells.Item(344, 2).Value = 'Owari no Seraph'
$ws.Cells.Item(344, 3).Value = 1
$ws.Cells.Item(344, 4).Value = 2
$ws.Cells.Item(344, 5).Value = 'Spring 2015'

$ws.Cells.Item(345, 1).Value = 343
$ws.Cells.Item(345, 2).Value = 'Danna ga Nani wo Itteiru ka Wakaranai Ken 2'
$ws.Cells.Item(345, 3).Value = 1
$ws.Cells.Item(345, 4).Value = 2
$ws.Cells.Item(345, 5).Value = 'Spring 2015'

$ws.Cells.Item(346, 1).Value = 344
$ws.Cells.Item(346, 2).Value = 'Teekyuu 4'
$ws.Cells.Item(346, 3).Value = 1
$ws.Cells.Item(346, 4).Value = 2
$ws.Cells.Item(346, 5).Value = 'Spring 2015'

$ws.Cells.Item(347, 1).Value = 345
$ws.Cells.Item(347, 2).Value = 'Kyoukai no Kanata Movie: I''ll Be Here - Kako-hen'
$ws.Cells.Item(347, 3).Value = 1
$ws.Cells.Item(347, 4).Value = 2
$ws.Cells.Item(347, 5).Value = 'Spring 2015'

$ws.Cells.Item(348, 1).Value = 346
$ws.Cells.Item(348, 2).Value = 'Kyoukai no Kanata Movie: I''ll Be Here - Mirai-hen'
$ws.Cells.Item(348, 3).Value = 1
$ws.Cells.Item(348, 4).Value = 2
$ws.Cells.Item(348, 5).Value = 'Spring 2015'

$ws.Cells.Item(349, 1).Value = 347
$ws.Cells.Item(349, 2).Value = 'Love Live! The School Idol Movie'
$ws.Cells.Item(349, 3).Value = 1
$ws.Cells.Item(349, 4).Value = 2
$ws.Cells.Item(349, 5).Value = 'Spring 2015'

$ws.Cells.Item(350, 1).Value = 348
$ws.Cells.Item(350, 2).Value = 'Persona 3 the Movie 3: Falling Down'
$ws.Cells.Item(350, 3).Value = 1
$ws.Cells.Item(350, 4).Value = 2
$ws.Cells.Item(350, 5).Value = 'Spring 2015'

$ws.Cells.Item(351, 1).Value = 349
$ws.Cells.Item(351, 2).Value = 'Tamayura: Sotsugyou Shashin Part 1 - Me: Kizashi'
$ws.Cells.Item(351, 3).Value = 1
$ws.Cells.Item(351, 4).Value = 2
$ws.Cells.Item(351, 5).Value = 'Spring 2015'

$ws.Cells.Item(352, 1).Value = 350
$ws.Cells.Item(352, 2).Value = 'Fate/stay night: Unlimited Blade Works 2nd Season'
$ws.Cells.Item(352, 3).Value = 1
$ws.Cells.Item(352, 4).Value = 1
$ws.Cells.Item(352, 5).Value = 'Spring 2015'
